$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-03 Wednesday" "2024-07-04 Thursday"

Replace-Text "89×86=" "75×78="
Replace-Text "56×79=" "41×11="
Replace-Text "78×99=" "95×16="
Replace-Text "65×42=" "54×57="
Replace-Text "64×71=" "91×42="

Replace-Text "17×26=" "61×28="
Replace-Text "48×11=" "46×43="
Replace-Text "24×39=" "67×57="
Replace-Text "29×41=" "15×18="
Replace-Text "59×89=" "32×35="

Replace-Text "83×19=" "66×24="
Replace-Text "75×72=" "29×39="
Replace-Text "52×46=" "62×63="
Replace-Text "25×28=" "56×89="
Replace-Text "85×98=" "37×15="

Replace-Text "60×18=" "65×43="
Replace-Text "70×34=" "30×70="
Replace-Text "89×58=" "94×73="
Replace-Text "56×85=" "69×14="
Replace-Text "74×23=" "20×45="

Replace-Text "30×71=" "77×56="
Replace-Text "60×40=" "18×85="
Replace-Text "34×25=" "52×44="
Replace-Text "11×70=" "50×67="
Replace-Text "88×67=" "90×94="
